$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11:17 down to 12:18
$ws.Rows.Item(11).Insert()

# Fill the new row 11 with data (same dimension values as the surrounding
# records, with the new date/price observations for this week)
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = 'Macroferia Regional de Talca'
$ws.Range("C11").Value = 'Maule'
$ws.Range("D11").Value = 44463
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = 'Espárragos'
$ws.Range("H11").Value = 'Verde'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 2000
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = 'Provincia de Linares'
$ws.Range("P11").Value = 2000
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 'Hortaliza'
